$wb = $excel.ActiveWorkbook

# --- Worksheets ---
$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status / datetime text updates (Ready for handoff, refreshed timestamps) ---

# Overview sheet: zh-cn / de-de status columns (E2, F2) and handoff generate date (G2)
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-23 18:42:10"

# zh-cn detail sheet: Status (C2) + Latest Handoff Datetime (H2)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-23 18:41:57"

# de-de detail sheet: Status (C2) + Latest Handoff Datetime (H2)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-23 18:42:10"

# --- Column width updates (widened to fit the new "Ready for handoff" status text) ---
# Target stored width ~17.216 characters; Excel's ColumnWidth COM property is
# quantized to a 1/6-character pixel grid, so 16.333... is the closest input
# that rounds to that stored width.
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
